$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Meetups": append " Meetup" to each week label in column D,
# and fix a typo in B17 ("8:30pm" -> "8:30 pm").
# ---------------------------------------------------------------------
$wsMeetups = $wb.Worksheets.Item("Meetups")

$wsMeetups.Range("D2").Value  = "Week 1 Meetup"
$wsMeetups.Range("D3").Value  = "Week 2 Meetup"
$wsMeetups.Range("D4").Value  = "Week 3 Meetup"
$wsMeetups.Range("D5").Value  = "Week 4 Meetup"
$wsMeetups.Range("D6").Value  = "Week 5 Meetup"
$wsMeetups.Range("D7").Value  = "Week 6 Meetup"
$wsMeetups.Range("D8").Value  = "Week 7 Meetup"
$wsMeetups.Range("D9").Value  = "Week 8 Meetup"
$wsMeetups.Range("D10").Value = "NO CLASS - Spring Recess"
$wsMeetups.Range("D11").Value = "Week 9 Meetup"
$wsMeetups.Range("D12").Value = "Week 10 Meetup"
$wsMeetups.Range("D13").Value = "Week 11 Meetup"
$wsMeetups.Range("D14").Value = "Week 12 Meetup"
$wsMeetups.Range("D15").Value = "Week 13 Meetup"
$wsMeetups.Range("D16").Value = "Week 14 Meetup"
$wsMeetups.Range("D17").Value = "Week 15 Meetup"

$wsMeetups.Range("B17").Value = "8:30 pm"

# ---------------------------------------------------------------------
# Sheet "Schedule": the "Topic" column used to hold Markdown-style
# links, e.g. "[Chapter 1](/chapters/chapter1)". Split those into a
# plain-text Topic column (C) plus a new "Link" column (D) holding just
# the URL, to be consumed by the new JS calendar. Also fix a couple of
# bad dates.
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("D1").Value = "Link"

$wsSchedule.Range("C2").Value = "Chapter 1"
$wsSchedule.Range("D2").Value = "/chapters/chapter1"

$wsSchedule.Range("A3").Value = 44235
$wsSchedule.Range("C3").Value = "Chapter 2"
$wsSchedule.Range("D3").Value = "/chapters/chapter2"

$wsSchedule.Range("C5").Value = "Chapter 3"
$wsSchedule.Range("D5").Value = "/chapters/chapter3"

$wsSchedule.Range("C6").Value = "Chapter 4"

$wsSchedule.Range("C7").Value = "Chapter 5"

$wsSchedule.Range("B10").Value = 44283

$wsSchedule.Range("C19").Value = "Final Exam"
$wsSchedule.Range("D19").Value = "/assignments/exams/"

$wsSchedule.Columns.Item(3).ColumnWidth = 32.66666667

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping to match the edited file: leave
# "Schedule" with B2 selected, then finish on "Meetups" (the originally
# active sheet) with D18 selected so it stays the active tab.
# ---------------------------------------------------------------------
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("B2").Select() | Out-Null

$wsMeetups.Activate() | Out-Null
$wsMeetups.Range("D18").Select() | Out-Null
